# Applies:
#  1) Three tables (slides 14-16) switch from the custom "Table_0" style
#     ({FEB0C677-119A-4F44-9A3E-95C50E4BC091}) to the built-in table
#     style {91595EB4-651B-422A-897B-F5B8303502A9}.
#  2) The deck's theme (slide master theme, ppt/theme/theme1.xml) changes
#     its colour scheme from the "Integral" (Red Violet) palette to the
#     default "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newTableStyle = "{91595EB4-651B-422A-897B-F5B8303502A9}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Theme colour scheme -------------------------------------------
# Office Theme default palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
function RGBVal([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$officeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1
    (RGBVal 0xFF 0xFF 0xFF),  # lt1
    (RGBVal 0x44 0x54 0x6A),  # dk2
    (RGBVal 0xE7 0xE6 0xE6),  # lt2
    (RGBVal 0x5B 0x9B 0xD5),  # accent1
    (RGBVal 0xED 0x7D 0x31),  # accent2
    (RGBVal 0xA5 0xA5 0xA5),  # accent3
    (RGBVal 0xFF 0xC0 0x00),  # accent4
    (RGBVal 0x44 0x72 0xC4),  # accent5
    (RGBVal 0x70 0xAD 0x47),  # accent6
    (RGBVal 0x05 0x63 0xC1),  # hlink
    (RGBVal 0x95 0x4F 0x72)   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
